$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency data (prices, 1h volume %, and three
# rows whose coin/link/price/volume were fully replaced because the
# ranking order shifted) as produced by the scheduled data refresh.

# Row 2
$ws.Range("D2").Value = "'61.257.32"
$ws.Range("E2").Value = "'  -0.87%  "

# Row 3
$ws.Range("D3").Value = "'2.394.51"
$ws.Range("E3").Value = "'  -4.14%  "

# Row 4
$ws.Range("E4").Value = "'  +0.02%  "

# Row 5
$ws.Range("D5").Value = "'548.68"
$ws.Range("E5").Value = "'  -1.13%  "

# Row 6
$ws.Range("D6").Value = "'142.26"
$ws.Range("E6").Value = "'  -3.29%  "

# Row 7
$ws.Range("E7").Value = "'  -0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.543"
$ws.Range("E8").Value = "'  -10.17%  "

# Row 9
$ws.Range("D9").Value = "'2.393.94"
$ws.Range("E9").Value = "'  -4.08%  "

# Row 10
$ws.Range("E10").Value = "'  -1.84%  "

# Row 11
$ws.Range("E11").Value = "'  +0.34%  "

# Row 12
$ws.Range("D12").Value = "'5.27"
$ws.Range("E12").Value = "'  -3.60%  "

# Row 13
$ws.Range("E13").Value = "'  -3.05%  "

# Row 14
$ws.Range("E14").Value = "'  -3.07%  "

# Row 15
$ws.Range("D15").Value = "'2.825.63"
$ws.Range("E15").Value = "'  -4.04%  "

# Row 16
$ws.Range("E16").Value = "'  -1.27%  "

# Row 17
$ws.Range("D17").Value = "'61.066.08"
$ws.Range("E17").Value = "'  -1.03%  "

# Row 18
$ws.Range("D18").Value = "'2.388.72"
$ws.Range("E18").Value = "'  -3.86%  "

# Row 19
$ws.Range("D19").Value = "'10.76"
$ws.Range("E19").Value = "'  -4.21%  "

# Row 20
$ws.Range("E20").Value = "'  -2.02%  "

# Row 21
$ws.Range("D21").Value = "'319.72"
$ws.Range("E21").Value = "'  -1.33%  "

# Row 22
$ws.Range("D22").Value = "'6.78"
$ws.Range("E22").Value = "'  -3.69%  "

# Row 23
$ws.Range("D23").Value = "'1.95"
$ws.Range("E23").Value = "'  +7.64%  "

# Row 24
$ws.Range("E24").Value = "'  +0.04%  "

# Row 25
$ws.Range("D25").Value = "'63.88"
$ws.Range("E25").Value = "'  -0.46%  "

# Row 26
$ws.Range("D26").Value = "'8.13"
$ws.Range("E26").Value = "'  +6.23%  "

# Row 27
$ws.Range("B27").Value = "'Bittensor"
$ws.Range("C27").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D27").Value = "'543.16"
$ws.Range("E27").Value = "'  +0.60%  "

# Row 28
$ws.Range("B28").Value = "'PEPE"
$ws.Range("C28").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "'0.0₃0949"
$ws.Range("E28").Value = "'  -5.42%  "

# Row 29
$ws.Range("B29").Value = "'Binance-PegBSC-USD"
$ws.Range("C29").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "'  +0.11%  "

# Row 30
$ws.Range("B30").Value = "'WrappedeETH"
$ws.Range("C30").Value = "'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D30").Value = "'2.506.82"
$ws.Range("E30").Value = "'  -3.68%  "

# Row 31
$ws.Range("D31").Value = "'1.45"
$ws.Range("E31").Value = "'  -5.18%  "

# Row 32
$ws.Range("D32").Value = "'8.15"
$ws.Range("E32").Value = "'  -4.04%  "

# Row 33
$ws.Range("E33").Value = "'  -3.82%  "

# Row 34
$ws.Range("E34").Value = "'  -3.96%  "

# Row 35
$ws.Range("E35").Value = "'  -0.96%  "

# Row 36
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "'  +0.02%  "

# Row 37
$ws.Range("E37").Value = "'  -7.05%  "

# Row 38
$ws.Range("E38").Value = "'  -4.33%  "

# Row 39
$ws.Range("D39").Value = "'0.379"
$ws.Range("E39").Value = "'  -2.07%  "

# Row 40
$ws.Range("D40").Value = "'1.84"
$ws.Range("E40").Value = "'  +5.37%  "

# Row 41
$ws.Range("D41").Value = "'18.15"
$ws.Range("E41").Value = "'  -2.54%  "

# Row 42
$ws.Range("B42").Value = "'Monero"
$ws.Range("C42").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'139.20"
$ws.Range("E42").Value = "'  -6.43%  "

# Row 43
$ws.Range("B43").Value = "'USDe"
$ws.Range("C43").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "'  +0.05%  "

# Row 44
$ws.Range("D44").Value = "'40.38"
$ws.Range("E44").Value = "'  -0.18%  "

# Row 45
$ws.Range("D45").Value = "'2.25"
$ws.Range("E45").Value = "'  -5.95%  "

# Row 46
$ws.Range("D46").Value = "'142.44"
$ws.Range("E46").Value = "'  -4.14%  "

# Row 47
$ws.Range("D47").Value = "'3.61"
$ws.Range("E47").Value = "'  -0.82%  "

# Row 48
$ws.Range("D48").Value = "'20.31"
$ws.Range("E48").Value = "'  -2.99%  "

# Row 49
$ws.Range("D49").Value = "'0.0522"
$ws.Range("E49").Value = "'  -2.71%  "

# Row 50
$ws.Range("D50").Value = "'0.579"
$ws.Range("E50").Value = "'  -3.67%  "

# Row 51
$ws.Range("B51").Value = "'VeChain"
$ws.Range("C51").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0228"
$ws.Range("E51").Value = "'  -0.90%  "
